# Fruta / hortaliza, semanal
# A new week of data (row 8) is prepended; every existing data row (8..124)
# shifts down by one, and the last row that falls off the bottom becomes a
# brand-new row 125.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstDataRow = 8
$lastDataRow = 124

# --- capture the "before" state of the columns that move with each record:
#     D = Fecha, J = Volumen, K = Precio minimo, L = Precio maximo,
#     M = Precio promedio ponderado, P = Precio $/Kg
$colD = @()
$colJ = @()
$colK = @()
$colL = @()
$colM = @()
$colP = @()

for ($r = $firstDataRow; $r -le $lastDataRow; $r++) {
    $colD += $ws.Cells.Item($r, 4).Value2
    $colJ += $ws.Cells.Item($r, 10).Value2
    $colK += $ws.Cells.Item($r, 11).Value2
    $colL += $ws.Cells.Item($r, 12).Value2
    $colM += $ws.Cells.Item($r, 13).Value2
    $colP += $ws.Cells.Item($r, 16).Value2
}

# --- the brand-new first record (new week of prices) ---
$newDate = 44616
$newJ = 2200
$newK = 2300
$newL = 2500
$newM = 2400
$newP = 1600

# --- shift every existing record down by one row; the record that was in
#     the last data row becomes the new row ($lastDataRow + 1) ---
$newLastRow = $lastDataRow + 1

for ($r = $newLastRow; $r -ge ($firstDataRow + 1); $r--) {
    $srcIdx = ($r - 1) - $firstDataRow

    $ws.Cells.Item($r, 4).Value2 = $colD[$srcIdx]
    $ws.Cells.Item($r, 10).Value2 = $colJ[$srcIdx]
    $ws.Cells.Item($r, 11).Value2 = $colK[$srcIdx]
    $ws.Cells.Item($r, 12).Value2 = $colL[$srcIdx]
    $ws.Cells.Item($r, 13).Value2 = $colM[$srcIdx]
    $ws.Cells.Item($r, 16).Value2 = $colP[$srcIdx]
}

# the new row needs the rest of its (constant, identical-for-every-row)
# columns filled in too, since it did not exist before
$ws.Cells.Item($newLastRow, 1).Value2 = $ws.Cells.Item($firstDataRow, 1).Value2
$ws.Cells.Item($newLastRow, 2).Value2 = $ws.Cells.Item($firstDataRow, 2).Value2
$ws.Cells.Item($newLastRow, 3).Value2 = $ws.Cells.Item($firstDataRow, 3).Value2
$ws.Cells.Item($newLastRow, 5).Value2 = $ws.Cells.Item($firstDataRow, 5).Value2
$ws.Cells.Item($newLastRow, 6).Value2 = $ws.Cells.Item($firstDataRow, 6).Value2
$ws.Cells.Item($newLastRow, 7).Value2 = $ws.Cells.Item($firstDataRow, 7).Value2
$ws.Cells.Item($newLastRow, 8).Value2 = $ws.Cells.Item($firstDataRow, 8).Value2
$ws.Cells.Item($newLastRow, 9).Value2 = $ws.Cells.Item($firstDataRow, 9).Value2
$ws.Cells.Item($newLastRow, 14).Value2 = $ws.Cells.Item($firstDataRow, 14).Value2
$ws.Cells.Item($newLastRow, 15).Value2 = $ws.Cells.Item($firstDataRow, 15).Value2
$ws.Cells.Item($newLastRow, 17).Value2 = $ws.Cells.Item($firstDataRow, 17).Value2
$ws.Cells.Item($newLastRow, 18).Value2 = $ws.Cells.Item($firstDataRow, 18).Value2

# give the new row's date cell the same date number-format as the rest of
# the Fecha column
$ws.Cells.Item($newLastRow, 4).NumberFormat = $ws.Cells.Item($firstDataRow, 4).NumberFormat

# --- finally, write the brand-new first record into row 8 ---
$ws.Cells.Item($firstDataRow, 4).Value2 = $newDate
$ws.Cells.Item($firstDataRow, 10).Value2 = $newJ
$ws.Cells.Item($firstDataRow, 11).Value2 = $newK
$ws.Cells.Item($firstDataRow, 12).Value2 = $newL
$ws.Cells.Item($firstDataRow, 13).Value2 = $newM
$ws.Cells.Item($firstDataRow, 16).Value2 = $newP
